# edit.ps1
# Adds a new "2021" column (N) of data to the sheet, mirroring the layout of
# the existing 2011-2020 (D-M) columns, changes the number formatting style
# used by the D4:L4 data row, and updates the page setup / selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats

# ---------------------------------------------------------------------------
# 1) Row 4 (D4:L4): switch the number style used for the yearly figures to
#    the bold "total" style already used by M4 (style index 15), without
#    touching the A4:C4 label cells.
# ---------------------------------------------------------------------------
$ws.Range("D4:L4").Font.Bold = $true

# ---------------------------------------------------------------------------
# 2) Add the new column N (year 2021) for rows 2-15, re-using the same cell
#    formatting as the corresponding cell in the row (copied via
#    PasteSpecial so the existing style is reused rather than a new one
#    being created).
# ---------------------------------------------------------------------------

# Row 2 : blank separator cell, same style as M2
$ws.Range("N2").Value = $null
$ws.Range("M2").Copy()
$ws.Range("N2").PasteSpecial($xlPasteFormats)

# Row 3 : year header
$ws.Range("N3").Value = 2021
$ws.Range("M3").Copy()
$ws.Range("N3").PasteSpecial($xlPasteFormats)

# Row 4 : total row, bold number style (same as M4)
$ws.Range("N4").Value = 95.134712433469176
$ws.Range("M4").Copy()
$ws.Range("N4").PasteSpecial($xlPasteFormats)

# Row 5
$ws.Range("N5").Value = 99.705541665880986
$ws.Range("M5").Copy()
$ws.Range("N5").PasteSpecial($xlPasteFormats)

# Row 6
$ws.Range("N6").Value = 92.425193326577897
$ws.Range("M6").Copy()
$ws.Range("N6").PasteSpecial($xlPasteFormats)

# Row 7 : note style matches D7 (plain), not M7 (right aligned)
$ws.Range("N7").Value = 88.209991167538519
$ws.Range("D7").Copy()
$ws.Range("N7").PasteSpecial($xlPasteFormats)

# Row 8
$ws.Range("N8").Value = 92.225038985690773
$ws.Range("D8").Copy()
$ws.Range("N8").PasteSpecial($xlPasteFormats)

# Row 9
$ws.Range("N9").Value = 96.801032063987265
$ws.Range("D9").Copy()
$ws.Range("N9").PasteSpecial($xlPasteFormats)

# Row 10
$ws.Range("N10").Value = 97.660491031729507
$ws.Range("D10").Copy()
$ws.Range("N10").PasteSpecial($xlPasteFormats)

# Row 11
$ws.Range("N11").Value = 90.23262877800066
$ws.Range("D11").Copy()
$ws.Range("N11").PasteSpecial($xlPasteFormats)

# Row 12
$ws.Range("N12").Value = 99.653994395099105
$ws.Range("D12").Copy()
$ws.Range("N12").PasteSpecial($xlPasteFormats)

# Row 13
$ws.Range("N13").Value = 100
$ws.Range("D13").Copy()
$ws.Range("N13").PasteSpecial($xlPasteFormats)

# Row 14
$ws.Range("N14").Value = 100
$ws.Range("D14").Copy()
$ws.Range("N14").PasteSpecial($xlPasteFormats)

# Row 15 : bottom border row, same style as M15
$ws.Range("N15").Value = 100
$ws.Range("M15").Copy()
$ws.Range("N15").PasteSpecial($xlPasteFormats)

# ---------------------------------------------------------------------------
# 3) Update the selected cell / page setup to match the saved view state.
# ---------------------------------------------------------------------------
$ws.Range("N2").Select() | Out-Null

$ws.PageSetup.HorizontalDpi = 300
$ws.PageSetup.VerticalDpi = 300
